$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dbo_unit")

# Insert a new row at row 2, shifting existing rows down.
$ws.Rows.Item(2).Insert()

# Populate the new row with the "Not applicable" lookup entry.
$ws.Cells.Item(2, 1).Value = -1
$ws.Cells.Item(2, 2).Value = "Not applicable"

# Update the named range to cover the new row.
$wb.Names.Item("dbo_unit").RefersTo = "=dbo_unit!`$A`$1:`$D`$16"

# Move the active selection to A3 (matches the saved selection state).
$ws.Range("A3").Select()
